$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 392 (shifts existing rows 392-428 down to 393-429),
# copying formatting (e.g. date style) from the row above, matching Excel's
# default Insert behavior.
$ws.Rows.Item(392).Insert()

# Populate the newly inserted row 392 with the new record.
$ws.Cells.Item(392, 1).Value = 4
$ws.Cells.Item(392, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(392, 3).Value = "Los Lagos"
$ws.Cells.Item(392, 4).Value = 44769
$ws.Cells.Item(392, 5).Value = 10
$ws.Cells.Item(392, 6).Value = 100114001
$ws.Cells.Item(392, 7).Value = "Papa"
$ws.Cells.Item(392, 8).Value = "Patagonia"
$ws.Cells.Item(392, 9).Value = "1a (guarda)"
$ws.Cells.Item(392, 10).Value = 150
$ws.Cells.Item(392, 11).Value = 8000
$ws.Cells.Item(392, 12).Value = 8000
$ws.Cells.Item(392, 13).Value = 8000
$ws.Cells.Item(392, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(392, 15).Value = "Provincia de Llanquihue"
$ws.Cells.Item(392, 16).Value = 320
$ws.Cells.Item(392, 17).Value = 25
$ws.Cells.Item(392, 18).Value = "Hortaliza"
